$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Base Load / Unavailable / Year): 2021 -> 2003 across M:W ---
$ws.Range("M3:W3").Value = 2003

# --- Row 4 (Base Load / Retrofit_existing_min / %): new ramp, drop the
#     trailing "=prior cell" formulas in favor of plain literal 1's ---
$ws.Range("N4").Value = 0.25
$ws.Range("O4").Value = 0.5
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 1

# --- Row 5 (Shoulder Load / Unavailable / Year): 2021 -> 2003 across M:W ---
$ws.Range("M5:W5").Value = 2003

# --- Row 6 (Shoulder Load / Retrofit_existing_min / %): same ramp/format
#     change as row 4 ---
$ws.Range("N6").Value = 0.25
$ws.Range("O6").Value = 0.5
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 1

# --- View state: scrolled/selected cell moves from A1:X6 to X4 ---
$ws.Range("X4").Select()
